$d = $word.ActiveDocument

# --- Simple text replacements (character-tag changes in dialogue lines) ---

$d.Content.Find.Execute(
    "Mom (waving smile): Good morning.", $true, $false, $false, $false, $false,
    $true, 1, $false, "Mom (neutral smiling): Good morning.", 2) | Out-Null

$d.Content.Find.Execute(
    "Mom (neutral smiling): Here’s your breakfast.", $true, $false, $false, $false, $false,
    $true, 1, $false, "Mom (neutral smiling_eyes_closed): Here’s your breakfast.", 2) | Out-Null

$d.Content.Find.Execute(
    "Mom (neutral smiling): After school today…", $true, $false, $false, $false, $false,
    $true, 1, $false, "Mom (neutral thinking): After school today…", 2) | Out-Null

$d.Content.Find.Execute(
    "Mom (neutral smiling): Will you be coming straight home?", $true, $false, $false, $false, $false,
    $true, 1, $false, "Mom (neutral curious): Will you be coming straight home?", 2) | Out-Null

# --- Merge runs that should become a single run (text content unchanged) ---

$d.Content.Find.Execute(
    "Mara (neutral confused): ...how you look?", $true, $false, $false, $false, $false,
    $true, 1, $false, "Mara (neutral confused): ...how you look?", 2) | Out-Null

$d.Content.Find.Execute(
    "I tell her about everything that happened yesterday, everything we talked about and said.", $true, $false, $false, $false, $false,
    $true, 1, $false, "I tell her about everything that happened yesterday, everything we talked about and said.", 2) | Out-Null

# --- Insert a new paragraph "Mom (exit):" right after "Mom (neutral smiling): Alright then." ---

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Mom (neutral smiling): Alright then.*") {
        $p.Range.InsertParagraphAfter()
        $p.Next().Range.Text = "Mom (exit):"
        break
    }
}
